$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert two new blank rows before the old row 13 ("Allergi med betydning...")
# This shifts the old rows 13-23 down to 15-25.
# ---------------------------------------------------------------------------
$ws.Rows("13:14").Insert()

# ---------------------------------------------------------------------------
# First lay down all the formatting (re-using existing format combinations so
# no spurious new border/fill entries get minted), using cells that already
# hold pre-existing shared strings so the shared-string table order below is
# governed solely by the Value2 assignments that follow.
# ---------------------------------------------------------------------------

# A13 / A14 - blank cells, same format as A7 (no fill, box border, font2)
$ws.Range("A7").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A7").Copy()
$ws.Range("A14").PasteSpecial(-4122)

# C13 / C14 / D13 / D14 - same format as C7 / D7
$ws.Range("C7").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D13").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("D14").PasteSpecial(-4122)

# E13 / E14 - same format as C7 but with a medium right-hand border (new style)
$ws.Range("C7").Copy()
$ws.Range("E13").PasteSpecial(-4122)
$ws.Range("E13").Borders.Item(10).LineStyle = 1
$ws.Range("E13").Borders.Item(10).Weight = -4138
$ws.Range("C7").Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("E14").Borders.Item(10).LineStyle = 1
$ws.Range("E14").Borders.Item(10).Weight = -4138

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Now fill in the cell text, in the exact order the new strings were added to
# the workbook, so the shared-string table comes out in the same order.
# ---------------------------------------------------------------------------
$ws.Range("C13").Value2 = "Antibiotikaprofylakse"
$ws.Range("C14").Value2 = "Tromboseprofylakse"
$ws.Range("B13").Value2 = "Medikamenter ved kirurgi"
$ws.Range("D13").Value2 = "Sjekk at lenken virker:)"
$ws.Range("D14").Value2 = "Sjekk at lenken virker:)"
$ws.Range("E13").Value2 = "MGR"
$ws.Range("E14").Value2 = "MGR"

# ---------------------------------------------------------------------------
# Update the active selection to match the edited workbook.
# ---------------------------------------------------------------------------
$ws.Range("E14").Select()
